# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-30 19:27:26
#
# The "Recorded By" column (G) lists the same set of recorders but in a
# different, inconsistent order across rows. Normalize the ordering of the
# comma-separated names/emails within each G-column cell:
#   "dnasr281@gmail.com, System"                  -> "System, dnasr281@gmail.com"
#   "System, backup@backdoor.com, system"         -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val) { continue }

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
